$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (publishedDate) stays text, not auto-converted to a date/number
$ws.Range("D2:D31").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = '-Q7y6beGZHAC'
$ws.Range("B2").Value = 'Estudio del sector editorial de los juegos de rol en España :historia, tipología, perfil de lector, del autor, del traductor y del editor'
$ws.Range("C2").Value = 'Desconocido'
$ws.Range("D2").Value = '2009-06-18'
$ws.Range("E2").Value = 'Héctor Sevillano Pareja'

# Row 3
$ws.Range("A3").Value = '5XWXuQAACAAJ'
$ws.Range("B3").Value = 'El señor de los anillos'
$ws.Range("C3").Value = 'Cuatro volúmenes en uno.Incluye los mapas y apéndices originales de J.R.R. Tolkien, y un índice de canciones, nombres de personajes, lugares y cosas. 50 ilustraciones en color de Alan Lee.'
$ws.Range("D3").Value = '1993'
$ws.Range("E3").Value = 'J. R. R. Tolkien, Alan Lee'

# Row 4
$ws.Range("A4").Value = '6463zwEACAAJ'
$ws.Range("B4").Value = 'The return of the king'
$ws.Range("C4").Value = 'Desconocido'
$ws.Range("D4").Value = '1965'
$ws.Range("E4").Value = 'John Ronald Reuel Tolkien'

# Row 5
$ws.Range("A5").Value = '70jL1nnXg3wC'
$ws.Range("B5").Value = 'Daniel in the Lions'' Den'
$ws.Range("C5").Value = 'Desconocido'
$ws.Range("D5").Value = '1996'
$ws.Range("E5").Value = 'Ronne Randall'

# Row 6
$ws.Range("A6").Value = 'DYmUGGwZ8_oC'
$ws.Range("B6").Value = 'El Señor de los Anillos no 01/03 La Comunidad del Anillo (edición revisada)'
$ws.Range("C6").Value = 'Primera entrega de la trilogía. «Este libro es como un relámpago en un cielo claro. Decir que la novela heroica, espléndida, elocuente y desinhibida, ha retornado de pronto en una época de un antirromanticismo casi patológico, sería inadecuado. Para quienes vivimos en esa extraña época, el retorno —y el alivio que nos trae— es sin duda lo más importante. Pero para la historia misma de la novela —una historia que se remonta a la Odisea y a antes de la Odisea— no es un retorno, sino un paso adelante o una revolución: la conquista de un territorio nuevo.» —C.S. Lewis, Time & Tide, 1954 «La obra de Tolkien, difundida en millones de ejemplares, traducida a docenas de lenguas, inspiradora de slogans pintados en las paredes de Nueva York y de Buenos Aires... una coherente mitología de una autenticidad universal creada en pleno siglo veinte.» —George Steiner, Le Monde, 1973'
$ws.Range("D6").Value = '2010-07-15'
$ws.Range("E6").Value = 'J. R. R. Tolkien'

# Row 7
$ws.Range("A7").Value = 'FlGcUAnpMmIC'
$ws.Range("B7").Value = 'Hobbit''s Travels'
$ws.Range("C7").Value = 'Printed on deluxe recycled parchment paper, this journal celebrating J. R. R. Tolkien''s classic tales makes a lovely gift, and is just as nice to keep! With magical two-color illustrations throughout (drawings made by Frodo Baggins''s devoted companion, Sam Gamgee, on their travels throughout Middle-earth), it provides ample space for recording personal thoughts, reflections on Tolkien''s masterpiece, or fantasies of your own creation.'
$ws.Range("D7").Value = '2002-01-04'
$ws.Range("E7").Value = 'Sam Gamgee'

# Row 8
$ws.Range("A8").Value = 'Gdvd0AEACAAJ'
$ws.Range("B8").Value = 'El señor de los anillos: La comunidad del anillo'
$ws.Range("C8").Value = 'Desconocido'
$ws.Range("D8").Value = '1993'
$ws.Range("E8").Value = 'John Ronald Reuel Tolkien'

# Row 9
$ws.Range("A9").Value = 'IOyjMQEACAAJ'
$ws.Range("B9").Value = 'The Power of the Ring'
$ws.Range("C9").Value = 'Digging deep into J. R. R. Tolkien''s spiritual biography--his religious scholarship and his love of both Christian and pagan myth--Stratford Caldecott offers a critical study of how the acclaimed author effectively created a vivid Middle Earth using the familiar rites and ceremonies of human history. And while readers and moviegoers alike may appreciate the fantasy world of The Hobbit and the Lord of the Rings trilogy, few know that in life, Tolkien was a devout Roman Catholic and that the characters, the events, and the general morality of each novel are informed by the dogmas of his faith. Revised and updated, this acclaimed study of Tolkien''s achievement--previously released as Secret Fire in the UK--includes commentary on Peter Jackson''s film adaptations and explores many of the fascinating stories and letters published after Tolkien''s death.'
$ws.Range("D9").Value = '2012-12-04'
$ws.Range("E9").Value = 'Stratford Caldecott'

# Row 10
$ws.Range("A10").Value = 'JUdOAAAACAAJ'
$ws.Range("B10").Value = 'The Lord of the Rings Sketchbook'
$ws.Range("C10").Value = '"In The Lord of the Rings Sketchbook Alan Lee reveals in pictures and in words how he created the watercolor paintings for the special centenary edition of The Lord of the Rings. These images would prove so powerful and evocative that they would eventually define the look of Peter Jackson''s movie trilogy and earn him a coveted Academy Award." "The book is filled with more than 150 of his sketches and early conceptual pieces showing how the project progressed from idea to finished art. It also contains a selection of full-page paintings reproduced in full color, together with numerous examples of previously unseen conceptual art produced for the films and many new works drawn specially for this book." "The Lord of the Rings Sketchbook provides an insight into the imagination of the man who painted Tolkien''s vision, first on the page and then in three dimensions on the movie screen. It will also be of interest to many of the thousands of people who have bought the illustrated Lord of the Rings as well as to budding artists who want to unlock the secrets of book illustration."--BOOK JACKET.'
$ws.Range("D10").Value = '2005'
$ws.Range("E10").Value = 'Alan Lee'

# Row 11
$ws.Range("A11").Value = 'LCZvVRqH-m8C'
$ws.Range("B11").Value = 'El señor de los anillos'
$ws.Range("C11").Value = 'Desconocido'
$ws.Range("D11").Value = '1999-05'
$ws.Range("E11").Value = 'Terry Donaldson'

# Row 12
$ws.Range("A12").Value = 'LvsQ34A1fOMC'
$ws.Range("B12").Value = 'El Señor de los Anillos no 03/03 El Retorno del Rey (edición revisada)'
$ws.Range("C12").Value = 'La tercera entrega de la trilogía El Señor de los Anillos. Los ejércitos del Señor Oscuro van extendiendo cada vez más su maléfica sombra por la Tierra Media. Hombres, elfos y enanos unen sus fuerzas para presentar batalla a Sauron y sus huestes. Ajenos a estos preparativos, Frodo y Sam siguen adentrándose en el país de Mordor en su heroico viaje para destruir el Anillo de Poder en las Grietas del Destino. «Un final triunfante... un gran trabajo, tanto en la concepción como en la ejecución.» —Daily Telegraph «Un trabajo extraordinariamente imaginativo, parte saga, parte alegoría, y emocionante en su totalidad.» —The Times'
$ws.Range("D12").Value = '2010-07-15'
$ws.Range("E12").Value = 'J. R. R. Tolkien'

# Row 13
$ws.Range("A13").Value = 'Ndgf0AEACAAJ'
$ws.Range("B13").Value = 'El señor de los anillos'
$ws.Range("C13").Value = 'Desconocido'
$ws.Range("D13").Value = '2007'
$ws.Range("E13").Value = 'John Ronald Reuel Tolkien'

# Row 14
$ws.Range("A14").Value = 'QtSEvgEACAAJ'
$ws.Range("B14").Value = 'El señor de los anillos'
$ws.Range("C14").Value = 'Desconocido'
$ws.Range("D14").Value = '1998'
$ws.Range("E14").Value = 'John Ronald Reuel Tolkien'

# Row 15
$ws.Range("A15").Value = 'RYr8sgEACAAJ'
$ws.Range("B15").Value = 'Fellowship of the Ring'
$ws.Range("C15").Value = 'Desconocido'
$ws.Range("D15").Value = '2000'
$ws.Range("E15").Value = ''

# Row 16
$ws.Range("A16").Value = 'T8P3AAAACAAJ'
$ws.Range("B16").Value = 'Tolkien, el Señor de Los Anillos'
$ws.Range("C16").Value = 'Desconocido'
$ws.Range("D16").Value = '2004-09-01'
$ws.Range("E16").Value = 'J. R. R. Tolkien'

# Row 17
$ws.Range("A17").Value = 'WmdWtQAACAAJ'
$ws.Range("B17").Value = 'El Señor de los anillos'
$ws.Range("C17").Value = 'Desconocido'
$ws.Range("D17").Value = '2002'
$ws.Range("E17").Value = 'J. R. R. Tolkien'

# Row 18
$ws.Range("A18").Value = 'ZVwX0QEACAAJ'
$ws.Range("B18").Value = 'El Señor de los Anillos'
$ws.Range("C18").Value = 'Desconocido'
$ws.Range("D18").Value = '1985'
$ws.Range("E18").Value = 'J. R. R. Tolkien'

# Row 19
$ws.Range("A19").Value = 'ZcAlEAAAQBAJ'
$ws.Range("B19").Value = 'The Lord of the Rings Illustrated'
$ws.Range("C19").Value = 'A sumptuous slipcased edition of Tolkien''s classic epic tale of adventure, fully illustrated in color by the author himself. This deluxe volume is quarterbound in leather and includes many special features unique to this edition. Since it was first published in 1954, The Lord of the Rings has been a book people have treasured. Steeped in unrivalled magic and otherworldliness, its sweeping fantasy and epic adventure has touched the hearts of young and old alike. Over 100 million copies of its many editions have been sold around the world, and occasional collectors'' editions become prized and valuable items of publishing. This one-volume deluxe slipcased edition contains the complete text, fully corrected and reset, which is printed in red and black, and features thirty color illustrations, maps, and sketches drawn by Tolkien himself as he composed this epic work. These include the pages from the Book of Mazarbul, marvelous facsimiles created by Tolkien to accompany the famous "Bridge of Khazad-dum" chapter. Also appearing are two poster-size, fold-out maps revealing all the detail of Middle-earth. This very special deluxe edition is quarterbound in cloth and red leather, with raised ribs on the spine and stamped in two foils. The pages are edged in gold and contained within are special features unique to this edition.'
$ws.Range("D19").Value = '2021-10-19'
$ws.Range("E19").Value = 'J. R. R. Tolkien'

# Row 20
$ws.Range("A20").Value = 'aWZzLPhY4o0C'
$ws.Range("B20").Value = 'The Fellowship Of The Ring'
$ws.Range("C20").Value = 'Begin your journey into Middle-earth... The inspiration for the upcoming original series on Prime Video, The Lord of the Rings: The Rings of Power. The Fellowship of the Ring is the first part of J.R.R. Tolkien’s epic adventure The Lord of the Rings. One Ring to rule them all, One Ring to find them, One Ring to bring them all and in the darkness bind them. Sauron, the Dark Lord, has gathered to him all the Rings of Power—the means by which he intends to rule Middle-earth. All he lacks in his plans for dominion is the One Ring—the ring that rules them all—which has fallen into the hands of the hobbit, Bilbo Baggins. In a sleepy village in the Shire, young Frodo Baggins finds himself faced with an immense task, as his elderly cousin Bilbo entrusts the Ring to his care. Frodo must leave his home and make a perilous journey across Middle-earth to the Cracks of Doom, there to destroy the Ring and foil the Dark Lord in his evil purpose.'
$ws.Range("D20").Value = '2012-02-15'
$ws.Range("E20").Value = 'J.R.R. Tolkien'

# Row 21
$ws.Range("A21").Value = 'ayczzwEACAAJ'
$ws.Range("B21").Value = 'El señor de los anillos'
$ws.Range("C21").Value = 'Desconocido'
$ws.Range("D21").Value = '1988'
$ws.Range("E21").Value = 'John Ronald Ruelen Tolkien'

# Row 22
$ws.Range("A22").Value = 'bjUh2d_6atUC'
$ws.Range("B22").Value = 'Más allá de los niños índigo'
$ws.Range("C22").Value = 'Según la profecía, el quinto sol o quinto mundo del calendario maya entra en ascensión el 21 de diciembre de 2012. Esta fecha representa una “puerta de acceso” en el desarrollo planetario que expondrá a la humanidad a nuevas formas de vivir y nuevos mundos de oportunidades. Las predicciones de la antigüedad indican que nuestra transición satisfactoria a través de esta puerta de acceso depende de la “quinta raza raíz” (la nueva expresión de la reserva genética humana) destinada a ayudarnos a sobrellevar los enormes y emocionantes cambios que nos aguardan. En Más allá de los niños índigo, P. M. H. Atwater arroja luz sobre las características de los extraordinarios “nuevos niños” de la quinta raza raíz, esos chicos brillantes e irreverentes que han nacido desde 1982. La autora explora la relación de los nuevos niños con las profecías del calendario maya y otras tradiciones, proveyendo amplia información de antecedentes acerca de las siete razas raíz (la sexta y séptima de las cuales aún no han surgido) y la gran transformación de la conciencia que ya está teniendo lugar. La autora revela la conexión existente entre las siete razas raíz y los siete chakras, y la manera en que la humanidad verá abrirse el quinto chakra (el de la voluntad) a medida que los nuevos niños lleguen a la madurez. Analiza además el fenómeno de la inteligencia creciente y el potencial sin desarrollar y ofrece orientaciones y herramientas concretas para los que procuran comprender y ayudar a los nuevos niños a realizar todo su potencial. Más allá de los niños índigo es el primer estudio importante de los niños de hoy y del lugar que ocupan en nuestro mundo rápidamente cambiante; combina la investigación objetiva con la revelación mística y las profecías.'
$ws.Range("D22").Value = '2008-01-07'
$ws.Range("E22").Value = 'P. M. H. Atwater'

# Row 23
$ws.Range("A23").Value = 'cURzPgAACAAJ'
$ws.Range("B23").Value = 'El señor de los anillos'
$ws.Range("C23").Value = 'Desconocido'
$ws.Range("D23").Value = '1993'
$ws.Range("E23").Value = 'John Ronald Reuel Tolkien'

# Row 24
$ws.Range("A24").Value = 'e1ZJzwEACAAJ'
$ws.Range("B24").Value = 'El Señor de Los Anillos 3. El Retorno del Rey (TV Tie-In). the Lord of the Rings 3. the Return of the King (TV Tie-In) (Spanish Edition)'
$ws.Range("C24").Value = 'La última parte del viaje de Frodo y Sam Los ejércitos del Señor Oscuro van extendiendo cada vez más su maléfica sombra por la Tierra Media. Hombres, elfos y enanos unen sus fuerzas para presentar batalla a Sauron y sus huestes. Ajenos a estos preparativos, Frodo y Sam siguen adentrándose en el país de Mordor en su heroico viaje para destruir el Anillo de Poder en las Grietas del Destino. ENGLISH DESCRIPTION The Return of the King is the third part of J.R.R. Tolkien''s epic adventure The Lord of the Rings. One Ring to rule them all, One Ring to find them, One Ring to bring them all and in the darkness bind them. The Dark Lord has risen, and as he unleashes hordes of Orcs to conquer all Middle-earth, Frodo and Sam struggle deep into his realm in Mordor. To defeat Sauron, the One Ring must be destroyed in the fires of Mount Doom. But the way is impossibly hard, and Frodo is weakening. The Ring corrupts all who bear it and Frodo''s time is running out. Will Sam and Frodo succeed, or will the Dark Lord rule Middle-earth once more?'
$ws.Range("D24").Value = '2022-09-27'
$ws.Range("E24").Value = 'J. R. R. Tolkien'

# Row 25
$ws.Range("A25").Value = 'jZgjyAEACAAJ'
$ws.Range("B25").Value = 'The Return of the King'
$ws.Range("C25").Value = 'Desconocido'
$ws.Range("D25").Value = '1967'
$ws.Range("E25").Value = 'John Ronald Reuel Tolkien'

# Row 26
$ws.Range("A26").Value = 'ld5GswEACAAJ'
$ws.Range("B26").Value = 'The Fellowship of the Ring'
$ws.Range("C26").Value = 'Frodo the hobbit and a band of warriors from the different kingdoms set out to destroy the Ring of Power before the evil Sauron grasps control.'
$ws.Range("D26").Value = '2005'
$ws.Range("E26").Value = 'John Ronald Reuel Tolkien'

# Row 27
$ws.Range("A27").Value = 'neSkMQEACAAJ'
$ws.Range("B27").Value = 'El Señor de Los Anillos 1.'
$ws.Range("C27").Value = 'La Compania se ha disuelto y sus integrantes emprenden caminos separados. Frodo y Sam continuan solos su viaje a lo largo del rio Anduin, perseguidos por la sombra misteriosa de un ser extrano que tambien ambiciona la posesion del Anillo. Mientras, hombres, elfos y enanos se preparan para la batalla final contra las fuerzas del Senor del Mal.'
$ws.Range("D27").Value = '2012-11-13'
$ws.Range("E27").Value = 'John Ronald Reuel Tolkien'

# Row 28
$ws.Range("A28").Value = 'o5WfPwAACAAJ'
$ws.Range("B28").Value = 'El Señor de los anillos'
$ws.Range("C28").Value = 'Desconocido'
$ws.Range("D28").Value = '2002'
$ws.Range("E28").Value = 'J. R. R. Tolkien'

# Row 29
$ws.Range("A29").Value = 'q0JyPwAACAAJ'
$ws.Range("B29").Value = 'El Señor de Los Anillos, I'
$ws.Range("C29").Value = 'En la adormecida e idílica Comarca, un joven hobbit recibe un encargo: custodiar el Anillo Único y emprender el viaje para su destrucción en las Grietas del Destino. Acompañado por magos, hombres, elfos y enanos, atravesará la Tierra Media y se internará en las sombras de Mordor, perseguido siempre por las huestes de Sauron, el Señor Oscuro, dispuesto a recuperar su creación para establecer el dominio definitivo del Mal.'
$ws.Range("D29").Value = '2002'
$ws.Range("E29").Value = 'J. R. R. Tolkien'

# Row 30
$ws.Range("A30").Value = 'yK0Qi0R9U-IC'
$ws.Range("B30").Value = 'Impacto económico de las industrias culturales en Colombia'
$ws.Range("C30").Value = 'Economía y cultura - El sector cultural en tres escenarios - Sector editorial - Sector fonográfico - Cine - Sector televisión - La radio - Publicaciones periódicas - Recomendaciones de política sobre las industrias culturales en Colombia.'
$ws.Range("D30").Value = '2003'
$ws.Range("E30").Value = 'Colombia. Ministerio de Cultura'

# Row 31
$ws.Range("A31").Value = 'z5jCtAEACAAJ'
$ws.Range("B31").Value = 'El señor de los anillos'
$ws.Range("C31").Value = 'Desconocido'
$ws.Range("D31").Value = '1986'
$ws.Range("E31").Value = 'J. R. R. Tolkien'
